# Applies the edits described by the diff to the active document.
# Each call re-scopes the Find range to the whole document content so
# that every replacement searches from the start and only touches the
# single, uniquely-identified run it targets.

$d = $word.ActiveDocument

function Replace-Text($oldText, $newText, $wholeWord) {
    $rng = $d.Content
    $rng.Start = 0
    # Replace = 1 (wdReplaceOne) so only the single, first matching
    # occurrence is touched, even when MatchWholeWord would otherwise
    # also match a look-alike substring elsewhere (e.g. "01" inside a
    # "2025-09-01" date).
    $found = $rng.Find.Execute($oldText, $true, $wholeWord, $false, $false, $false, $true, 1, $false, $newText, 1)
    if (-not $found) {
        Write-Host "NOT FOUND: $oldText"
    }
}

# Codigo programa de formacion
Replace-Text "1000023" "3000006" $false

# Nombre del Programa
Replace-Text "APLICACION DE HERRAMIENTAS METODOLOGICAS EN INVEST" "EMPRENDEDOR EN PRESTACION DE SERVICIOS DE HERRAJE " $false

# Version del programa
Replace-Text "01" "03" $true

# Duracion Maxima (Horas)
Replace-Text "40" "300" $true

# Fecha de Inicio
Replace-Text "2025-09-01" "2025-09-06" $false

# Fecha prevista de terminacion
Replace-Text "2025-10-30" "2025-10-11" $false

# Departamento desarrollo de formacion
Replace-Text "Amazonas" "Cauca" $true

# Municipio desarrollo formacion
Replace-Text "Leticia" "Cajibío" $true

# Documento responsable label
Replace-Text "SIN DOCUMENTO #" "CC #" $false

# Documento responsable number
Replace-Text "123456" "234234234" $true

# Horario del curso de formacion
Replace-Text "12 a 8" "8 a 12" $false

# Fechas de ejecucion de la formacion (mes 2)
Replace-Text "22" "21" $true

# Firma Instructor trailing text
Replace-Text "asd" "123asd" $true

Write-Host "Done"
